$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.792.87"
$ws.Range("E2").Value = "'  -1.69%  "
$ws.Range("D3").Value = "'1.868.95"
$ws.Range("E3").Value = "'  -1.92%  "
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("D5").Value = "'300.01"
$ws.Range("E5").Value = "'  -2.53%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  -0.12%  "
$ws.Range("D7").Value = "'0.5343"
$ws.Range("E7").Value = "'  +1.56%  "
$ws.Range("D8").Value = "'0.3735"
$ws.Range("E8").Value = "'  -2.11%  "
$ws.Range("D9").Value = "'0.07137"
$ws.Range("E9").Value = "'  -2.08%  "
$ws.Range("D10").Value = "'21.53"
$ws.Range("E10").Value = "'  -2.36%  "
$ws.Range("D11").Value = "'0.8867"
$ws.Range("E11").Value = "'  -1.58%  "
$ws.Range("D12").Value = "'0.08154"
$ws.Range("E12").Value = "'  -0.53%  "
$ws.Range("D13").Value = "'1.889.66"
$ws.Range("E13").Value = "'  +38.74%  "
$ws.Range("D14").Value = "'92.37"
$ws.Range("E14").Value = "'  -3.52%  "
$ws.Range("D15").Value = "'5.287"
$ws.Range("E15").Value = "'  -1.17%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "'  -0.05%  "
$ws.Range("D17").Value = "'14.82"
$ws.Range("D18").Value = "'0.000008478"
$ws.Range("E18").Value = "'  -1.76%  "
$ws.Range("E19").Value = "'  -0.10%  "
$ws.Range("D20").Value = "'26.821.12"
$ws.Range("E20").Value = "'  -1.72%  "
$ws.Range("D21").Value = "'4.967"
$ws.Range("E21").Value = "'  -2.07%  "
$ws.Range("E22").Value = "'  -1.66%  "
$ws.Range("D23").Value = "'6.373"
$ws.Range("E23").Value = "'  -2.25%  "
$ws.Range("E24").Value = "'  -0.52%  "
$ws.Range("D25").Value = "'146.12"
$ws.Range("E25").Value = "'  -2.51%  "
$ws.Range("D26").Value = "'1.739"
$ws.Range("E26").Value = "'  +0.11%  "
$ws.Range("D27").Value = "'17.95"
$ws.Range("E27").Value = "'  -1.51%  "
$ws.Range("D28").Value = "'113.63"
$ws.Range("E28").Value = "'  -2.49%  "
$ws.Range("E29").Value = "'  -2.78%  "
$ws.Range("E30").Value = "'  -4.07%  "
$ws.Range("D31").Value = "'0.09086"
$ws.Range("E31").Value = "'  -1.80%  "
$ws.Range("E32").Value = "'  -2.99%  "
$ws.Range("D33").Value = "'0.05017"
$ws.Range("E33").Value = "'  -1.07%  "
$ws.Range("D34").Value = "'1.171"
$ws.Range("E34").Value = "'  -4.39%  "
$ws.Range("D35").Value = "'2.941"
$ws.Range("E35").Value = "'  -2.31%  "
$ws.Range("D36").Value = "'0.6091"
$ws.Range("E36").Value = "'  +5.96%  "
$ws.Range("D37").Value = "'2.677"
$ws.Range("E37").Value = "'  -0.10%  "
$ws.Range("D38").Value = "'3.188"
$ws.Range("E38").Value = "'  -4.93%  "
$ws.Range("D39").Value = "'0.01942"
$ws.Range("E39").Value = "'  -3.00%  "
$ws.Range("E40").Value = "'  -1.28%  "
$ws.Range("D41").Value = "'0.5276"
$ws.Range("E41").Value = "'  +7.53%  "
$ws.Range("D42").Value = "'8.762"
$ws.Range("E42").Value = "'  -5.70%  "
$ws.Range("D43").Value = "'6.474"
$ws.Range("E43").Value = "'  -0.89%  "
$ws.Range("D44").Value = "'116.23"
$ws.Range("E44").Value = "'  -0.37%  "
$ws.Range("D45").Value = "'0.1485"
$ws.Range("E45").Value = "'  -2.29%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "'  -0.12%  "
$ws.Range("D47").Value = "'1.645"
$ws.Range("E47").Value = "'  +0.67%  "
$ws.Range("D48").Value = "'9.944"
$ws.Range("E48").Value = "'  -2.34%  "
$ws.Range("D49").Value = "'37.29"
$ws.Range("E49").Value = "'  -3.88%  "
$ws.Range("D50").Value = "'0.06067"
$ws.Range("E50").Value = "'  -2.07%  "
$ws.Range("D51").Value = "'62.10"
$ws.Range("E51").Value = "'  -2.73%  "
